$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Working")

# Reference style cell (plain, unstyled text cell) used to restore the
# default "Normal" style after temporarily marking new cells as Text so
# that numeric/date-looking literals are stored as literal strings
# instead of being auto-converted to numbers / date serials.
$plainStyle = $ws.Range("A14").Style

# --- Row 15 ---------------------------------------------------------
$ws.Range("A15:E15").NumberFormat = "@"
$ws.Range("A15").Value = "Sarunas Stoncelis"
$ws.Range("B15").Value = "ref78999"
$ws.Range("C15").Value = "8794"
$ws.Range("D15").Value = "7685.3"
$ws.Range("E15").Value = "07/26/2022"
$ws.Range("A15:E15").Style = $plainStyle

# --- Row 16 ---------------------------------------------------------
$ws.Range("A16:E16").NumberFormat = "@"
$ws.Range("A16").Value = "Jonathan Fire"
$ws.Range("B16").Value = "paiment345"
$ws.Range("C16").Value = "0985"
$ws.Range("D16").Value = "658"
$ws.Range("E16").Value = "07/26/2022"
$ws.Range("A16:E16").Style = $plainStyle

# --- Row 17 (B/C/D left blank, only A and E populated) ---------------
$ws.Range("A17:E17").NumberFormat = "@"
$ws.Range("A17").Value = "jjoo"
$ws.Range("E17").Value = "07/26/2022"
$ws.Range("A17:E17").Style = $plainStyle
